# osi-and-tcp.pptx — "Add files via upload"
#
# 1) Table on slide 1: label the previously-blank "Layer Number" column
#    (col 3) with L7/L4/L3/L2/L1 values, merging the three Application/
#    Presentation/Session rows into one "L7(Firewall)" cell.
# 2) Refresh the cached date-placeholder text (datetimeFigureOut field)
#    on the slide master + every slide layout from 2020-07-12 to
#    2020-07-22.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Table edits (slide 1, shape 1 = "표 5")
# ---------------------------------------------------------------------
$s   = $p.Slides.Item(1)
$sh  = $s.Shapes.Item(1)
$tbl = $sh.Table

# Merge rows 2-4 of column 3 (Application/Presentation/Session) into a
# single cell, then label it "L7(Firewall)".
$topCell = $tbl.Cell(2, 3)
$botCell = $tbl.Cell(4, 3)
$topCell.Merge($botCell)
$tbl.Cell(2, 3).Shape.TextFrame.TextRange.Text = "L7(Firewall)"

# Remaining rows in column 3 just get their layer label.
$tbl.Cell(5, 3).Shape.TextFrame.TextRange.Text = "L4(NAT)"
$tbl.Cell(6, 3).Shape.TextFrame.TextRange.Text = "L3(router)"
$tbl.Cell(7, 3).Shape.TextFrame.TextRange.Text = "L2(switch)"
$tbl.Cell(8, 3).Shape.TextFrame.TextRange.Text = "L1"

# ---------------------------------------------------------------------
# 2. Cached "today" date placeholder: 2020-07-12 -> 2020-07-22
#    (slide master + all custom layouts)
# ---------------------------------------------------------------------
function Update-DatePlaceholder($container, $newText) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            $isDateShape = $false
            try {
                if ($shp.PlaceholderFormat.Type -eq 16) {
                    $isDateShape = $true
                }
            } catch {
            }
            if ($isDateShape) {
                $shp.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master "2020-07-22"

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout "2020-07-22"
}
